$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A slightly (was auto best-fit, now an explicit custom width)
$ws.Columns.Item(1).ColumnWidth = 9.8

# Row 7 - 20250408a trial
$ws.Range("A7").Value = "20250408a"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 6
$ws.Range("G7").Value = 7

# Row 9 - 20250408b trial
$ws.Range("A9").Value = "20250408b"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 6

# Row 10 - stray count value
$ws.Range("F10").Value = 10

# Row 11 - 20250409a trial
$ws.Range("A11").Value = "20250409a"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 6

# Row 12 - stray count values
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = 35

# Row 13 - 20250409b trial
$ws.Range("A13").Value = "20250409b"
$ws.Range("B13").Value = 8
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 9
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 11
$ws.Range("G13").Value = 12

# Row 14 - note
$ws.Range("G14").Value = "8,39"

# Row 15 - 20250424 trial (plain numeric date, not a shared string)
$ws.Range("A15").Value = 20250424
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 6

# Row 16 - stray count values
$ws.Range("C16").Value = 22
$ws.Range("D16").Value = 6

# Update selection to match final cursor position
[void]$ws.Range("B16").Select()
